$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.946.89'
$ws.Range("E2").Value = '  -0.05%  '
$ws.Range("D3").Value = '1.671.33'
$ws.Range("E3").Value = '  +1.15%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.80'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.06%  '
$ws.Range("E6").Value = '  +1.54%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("E8").Value = '  +0.40%  '
$ws.Range("E9").Value = '  +0.52%  '
$ws.Range("E10").Value = '  -0.05%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0891'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.61%  '
$ws.Range("D12").Value = '1.906.68'
$ws.Range("E12").Value = '  +1.17%  '
$ws.Range("D13").Value = '1.689.91'
$ws.Range("E13").Value = '  +2.24%  '
$ws.Range("E14").Value = '  +0.14%  '
$ws.Range("E15").Value = '  +1.07%  '
$ws.Range("E16").Value = '  +0.49%  '
$ws.Range("D17").Value = '26.945.19'
$ws.Range("E17").Value = '  -0.01%  '
$ws.Range("E18").Value = '  +3.98%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '234.02'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.95%  '
$ws.Range("E20").Value = '  +0.10%  '
$ws.Range("E21").Value = '  -0.03%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.43'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.16%  '
$ws.Range("E23").Value = '  -1.47%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.19'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.30%  '
$ws.Range("E25").Value = '  +0.68%  '
$ws.Range("E26").Value = '  +0.19%  '
$ws.Range("E27").Value = '  +0.73%  '
$ws.Range("E28").Value = '  -1.24%  '
$ws.Range("E29").Value = '  +0.14%  '
$ws.Range("E30").Value = '  +0.16%  '
$ws.Range("E31").Value = '  +0.26%  '
$ws.Range("E32").Value = '  +0.44%  '
$ws.Range("D33").Value = '1.460.28'
$ws.Range("E33").Value = '  -5.64%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.14'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.04%  '
$ws.Range("E35").Value = '  +1.85%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.41'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.17%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.581'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.09%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.898'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.14%  '
$ws.Range("E39").Value = '  +1.11%  '
$ws.Range("E40").Value = '  +11.91%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.80'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.18%  '
$ws.Range("E42").Value = '  +0.06%  '
$ws.Range("E43").Value = '  +2.91%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '66.47'
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = '1.812.64'
$ws.Range("E45").Value = '  +1.10%  '
$ws.Range("E46").Value = '  +0.91%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '90.59'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.93%  '
$ws.Range("E48").Value = '  +0.96%  '
$ws.Range("E49").Value = '  +2.77%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0508'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.53%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.66'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.73%  '
